$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 and 45 swap: PaxDollar <-> TrustWalletToken change rank order
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8425"
$ws.Range("E44").Value = "  +1.55%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D2").Value = "30.608.47"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.884.59"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.49"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4802"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06538"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.40"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "99.76"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07833"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7582"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "1.882.51"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.244"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.35"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "30.576.69"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.18"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007533"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "2.128.74"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.355"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.430"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.178"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.55"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.02"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.908"
$ws.Range("E28").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09782"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.326"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.505"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.254"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.181"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04850"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6984"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.768"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.870"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.309"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.53"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.971"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4249"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.36"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.030"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.30"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05789"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3962"
$ws.Range("E51").Value = "  +0.06%  "
